# Add two new member rows (row 3 and row 4) to Sheet1, matching the
# columns: Mem ID | NAME | LAST NAME | ID | NUMBER | DATE | history | balance

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @("436637", "morteza", "pashaei", "0441202453", "09961079096", "03/03/2023", "|", "0"),
    @("108171", "ali",     "reza",    "82998798",   "0923424145234", "03/03/2023", "|", "0")
)

$startRow = 3
$r = $startRow
foreach ($record in $newRows) {
    $c = 1
    foreach ($value in $record) {
        $cell = $ws.Cells.Item($r, $c)
        # Prefix with an apostrophe so values that look numeric/date-like
        # (leading zeros, slashes, etc.) are stored as literal text instead
        # of being auto-converted to numbers/dates, then strip the
        # quote-prefix formatting flag so no style gets attached to the
        # cell (keeps it identical to a plain, unstyled text cell).
        $cell.Value = "'" + $value
        $cell.ClearFormats()
        $c++
    }
    $r++
}
